$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 46

# Fill in the new work-diary entry first (so dependent formulas such as the
# SUM(C:C) total in H7 recalc against the final values).
$ws.Cells.Item($newRow, 1).Value = 44693
$ws.Cells.Item($newRow, 2).Value = "réalisation"
$ws.Cells.Item($newRow, 3).Value = 2.25
$ws.Cells.Item($newRow, 4).Value = "Paufinage de la vérification de l'email, activation de l'utilisateur, email tout bien"

# Copy formatting from the row above so the new row matches existing styles
$ws.Range("A44:D44").Copy()
$ws.Range("A46:D46").PasteSpecial(-4122) # xlPasteFormats

$ws.Rows.Item($newRow).RowHeight = 30

$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("A1:F46"))
